# Auto-generated edit script: updates market-price derived columns (H, I, J, K, L, M, N)
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets to reflect refreshed market data
# pulled by the scheduled runner.

$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H75").Value = 43768
$ws.Range("J75").Value = 43768
$ws.Range("L75").Value = 43768
$ws.Range("N75").Value = -45640
$ws.Range("H78").Value = 43768
$ws.Range("J78").Value = 43768
$ws.Range("L78").Value = 131304
$ws.Range("N78").Value = -140664
$ws.Range("H120").Value = 48997.25
$ws.Range("J120").Value = 48997.25
$ws.Range("L120").Value = 48997.25
$ws.Range("N120").Value = -58673.25
$ws.Range("H121").Value = 1010861
$ws.Range("I121").Value = 347.5
$ws.Range("J121").Value = 1076055.5
$ws.Range("K121").Value = 1042.5
$ws.Range("L121").Value = 3228166.5
$ws.Range("M121").Value = 704.5
$ws.Range("N121").Value = -3231660.5
$ws.Range("H135").Value = 17858326
$ws.Range("I135").Value = 1270.8572
$ws.Range("J135").Value = 71429490
$ws.Range("K135").Value = 11437.7148
$ws.Range("L135").Value = 642865410
$ws.Range("M135").Value = -8902.7148
$ws.Range("N135").Value = -642870480
$ws.Range("H137").Value = 2663.2292
$ws.Range("I137").Value = 1058.25
$ws.Range("J137").Value = 3198.2222
$ws.Range("K137").Value = 3174.75
$ws.Range("L137").Value = 9594.6666
$ws.Range("M137").Value = -624.75
$ws.Range("N137").Value = -14694.6666

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 25955.078
$ws.Range("I32").Value = 27038.072
$ws.Range("K32").Value = 27038.072
$ws.Range("M32").Value = -26751.072
$ws.Range("H74").Value = 2060.2593
$ws.Range("I74").Value = 1449.8667
$ws.Range("J74").Value = 2823.25
$ws.Range("K74").Value = 1449.8667
$ws.Range("L74").Value = 2823.25
$ws.Range("M74").Value = -575.8667
$ws.Range("N74").Value = -4571.25
$ws.Range("H77").Value = 2060.2593
$ws.Range("I77").Value = 1449.8667
$ws.Range("J77").Value = 2823.25
$ws.Range("K77").Value = 7249.333500000001
$ws.Range("L77").Value = 14116.25
$ws.Range("M77").Value = -2881.333500000001
$ws.Range("N77").Value = -22852.25

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H81").Value = 10835
$ws.Range("J81").Value = 10835
$ws.Range("L81").Value = 10835
$ws.Range("N81").Value = -12957
$ws.Range("H84").Value = 10835
$ws.Range("J84").Value = 10835
$ws.Range("L84").Value = 32505
$ws.Range("N84").Value = -43113
$ws.Range("H97").Value = 10903.375
$ws.Range("I97").Value = 1806.75
$ws.Range("K97").Value = 1806.75
$ws.Range("M97").Value = -815.75
$ws.Range("H130").Value = 48374
$ws.Range("J130").Value = 48374
$ws.Range("L130").Value = 48374
$ws.Range("N130").Value = -58414

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7884.607
$ws.Range("I31").Value = 3621.4167
$ws.Range("J31").Value = 11082
$ws.Range("K31").Value = 3621.4167
$ws.Range("L31").Value = 11082
$ws.Range("M31").Value = -3326.4167
$ws.Range("N31").Value = -11672
$ws.Range("H32").Value = 0
$ws.Range("I32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("H34").Value = 7884.607
$ws.Range("I34").Value = 3621.4167
$ws.Range("J34").Value = 11082
$ws.Range("K34").Value = 3621.4167
$ws.Range("L34").Value = 11082
$ws.Range("M34").Value = -3419.4167
$ws.Range("N34").Value = -11486
$ws.Range("H58").Value = 2125.6948
$ws.Range("I58").Value = 1758.1154
$ws.Range("J58").Value = 4856.2856
$ws.Range("K58").Value = 1758.1154
$ws.Range("L58").Value = 4856.2856
$ws.Range("M58").Value = -1555.1154
$ws.Range("N58").Value = -5262.2856
$ws.Range("H134").Value = 2842.3125
$ws.Range("I134").Value = 1207.4286
$ws.Range("K134").Value = 3622.2858
$ws.Range("M134").Value = -1087.2858
$ws.Range("H136").Value = 2125.6948
$ws.Range("I136").Value = 1758.1154
$ws.Range("J136").Value = 4856.2856
$ws.Range("K136").Value = 5274.3462
$ws.Range("L136").Value = 14568.8568
$ws.Range("M136").Value = -2724.3462
$ws.Range("N136").Value = -19668.8568
$ws.Range("M32").ClearContents()

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 818.875
$ws.Range("I5").Value = 417.46155
$ws.Range("J5").Value = 2558.3333
$ws.Range("K5").Value = 1252.38465
$ws.Range("L5").Value = 7674.999899999999
$ws.Range("M5").Value = -1140.38465
$ws.Range("N5").Value = -7898.999899999999
$ws.Range("H131").Value = 2216.3408
$ws.Range("I131").Value = 33908.668
$ws.Range("J131").Value = 1097.7882
$ws.Range("K131").Value = 101726.004
$ws.Range("L131").Value = 3293.3646
$ws.Range("M131").Value = -96686.00399999999
$ws.Range("N131").Value = -13373.3646
$ws.Range("H132").Value = 1298
$ws.Range("I132").Value = 1298
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 11682
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -9152
$ws.Range("H133").Value = 8205
$ws.Range("I133").Value = 8205
$ws.Range("J133").Value = 0
$ws.Range("K133").Value = 24615
$ws.Range("L133").Value = 0
$ws.Range("M133").Value = -19555
$ws.Range("H134").Value = 77158170
$ws.Range("I134").Value = 83587480
$ws.Range("J134").Value = 6500
$ws.Range("K134").Value = 250762440
$ws.Range("L134").Value = 19500
$ws.Range("M134").Value = -250757370
$ws.Range("N134").Value = -29640
$ws.Range("H135").Value = 818.875
$ws.Range("I135").Value = 417.46155
$ws.Range("J135").Value = 2558.3333
$ws.Range("K135").Value = 3757.15395
$ws.Range("L135").Value = 23024.9997
$ws.Range("M135").Value = -1222.15395
$ws.Range("N135").Value = -28094.9997
$ws.Range("H136").Value = 35717316
$ws.Range("I136").Value = 71430620
$ws.Range("J136").Value = 4013.8572
$ws.Range("K136").Value = 214291860
$ws.Range("L136").Value = 12041.5716
$ws.Range("M136").Value = -214286760
$ws.Range("N136").Value = -22241.5716
$ws.Range("H137").Value = 47627576
$ws.Range("I137").Value = 3210.8
$ws.Range("J137").Value = 90922456
$ws.Range("K137").Value = 9632.400000000001
$ws.Range("L137").Value = 272767368
$ws.Range("M137").Value = -4532.400000000001
$ws.Range("N137").Value = -272777568
$ws.Range("H139").Value = 8588.895
$ws.Range("I139").Value = 12471
$ws.Range("J139").Value = 3251
$ws.Range("K139").Value = 37413
$ws.Range("L139").Value = 9753
$ws.Range("M139").Value = -32273
$ws.Range("N139").Value = -20033
$ws.Range("H140").Value = 41755.6
$ws.Range("H141").Value = 50004060
$ws.Range("I141").Value = 71432230
$ws.Range("J141").Value = 4999.8335
$ws.Range("K141").Value = 214296690
$ws.Range("L141").Value = 14999.5005
$ws.Range("M141").Value = -214291510
$ws.Range("N141").Value = -25359.5005
$ws.Range("N132").ClearContents()
$ws.Range("N133").ClearContents()

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("H122").Value = 1146.6666
$ws.Range("I122").Value = 900
$ws.Range("J122").Value = 1196
$ws.Range("K122").Value = 2700
$ws.Range("L122").Value = 3588
$ws.Range("N122").Value = -8488
$ws.Range("M122").Value = -250
$ws.Range("H132").Value = 3053.0967
$ws.Range("I132").Value = 2327.389
$ws.Range("J132").Value = 4057.923
$ws.Range("K132").Value = 6982.167
$ws.Range("L132").Value = 12173.769
$ws.Range("M132").Value = -4452.167
$ws.Range("N132").Value = -17233.769
$ws.Range("N48").ClearContents()

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H59").Value = 14931.667
$ws.Range("J59").Value = 14931.667
$ws.Range("L59").Value = 14931.667
$ws.Range("N59").Value = -16239.667
$ws.Range("H92").Value = 40385
$ws.Range("J92").Value = 40385
$ws.Range("L92").Value = 40385
$ws.Range("N92").Value = -45377
$ws.Range("H102").Value = 48561
$ws.Range("J102").Value = 48561
$ws.Range("L102").Value = 48561
$ws.Range("N102").Value = -55051
$ws.Range("H129").Value = 37730.332
$ws.Range("J129").Value = 37730.332
$ws.Range("L129").Value = 37730.332
$ws.Range("N129").Value = -47730.332
$ws.Range("H131").Value = 32159
$ws.Range("J131").Value = 32159
$ws.Range("L131").Value = 32159
$ws.Range("N131").Value = -42239

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H93").Value = 34428.57
$ws.Range("J93").Value = 34428.57
$ws.Range("L93").Value = 34428.57
$ws.Range("N93").Value = -39420.57
$ws.Range("H102").Value = 40000
$ws.Range("J102").Value = 40000
$ws.Range("L102").Value = 40000
$ws.Range("N102").Value = -46490

Write-Output "Updated cells across ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets."